$d = $word.ActiveDocument

# 1) Fix "at least to replicas" -> "at least two replicas"
$d.Content.Find.Execute("at least to replicas, then ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "at least two replicas, then ", 2)

# 2) Bold "min-insync.replicas=2"
$rng = $d.Content
$rng.Find.Execute("min-insync.replicas=2", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $rng.Font.Bold = 1
}
